$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking text values (e.g. "296.76") stay as text,
# matching the source inlineStr cells, instead of being coerced to numbers
# by the default Range.Value assignment.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "39.342.43"
$ws.Range("E2").Value = "  -2.90%  "
$ws.Range("D3").Value = "2.203.14"
$ws.Range("E3").Value = "  -7.02%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "296.76"
$ws.Range("E5").Value = "  -4.63%  "
$ws.Range("D6").Value = "81.95"
$ws.Range("E6").Value = "  -4.50%  "
$ws.Range("D7").Value = "0.509"
$ws.Range("E7").Value = "  -4.51%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("D9").Value = "0.466"
$ws.Range("E9").Value = "  -4.86%  "
$ws.Range("D10").Value = "0.0769"
$ws.Range("E10").Value = "  -7.04%  "
$ws.Range("D11").Value = "28.95"
$ws.Range("E11").Value = "  -4.19%  "
$ws.Range("D12").Value = "47.36"
$ws.Range("E12").Value = "  -10.30%  "
$ws.Range("E13").Value = "  -2.64%  "
$ws.Range("D14").Value = "2.550.45"
$ws.Range("E14").Value = "  -6.94%  "
$ws.Range("D15").Value = "6.26"
$ws.Range("E15").Value = "  -3.34%  "
$ws.Range("D16").Value = "13.99"
$ws.Range("E16").Value = "  -5.94%  "
$ws.Range("D17").Value = "2.203.54"
$ws.Range("D18").Value = "0.710"
$ws.Range("E18").Value = "  -6.14%  "
$ws.Range("D19").Value = "39.224.14"
$ws.Range("E19").Value = "  -3.14%  "
$ws.Range("D20").Value = "0.0₃0871"
$ws.Range("E20").Value = "  -4.09%  "
$ws.Range("D21").Value = "5.70"
$ws.Range("E21").Value = "  -6.85%  "
$ws.Range("D22").Value = "64.88"
$ws.Range("E22").Value = "  -4.75%  "
$ws.Range("D23").Value = "10.28"
$ws.Range("E23").Value = "  -4.68%  "
$ws.Range("D24").Value = "225.38"
$ws.Range("E24").Value = "  -3.94%  "
$ws.Range("E25").Value = "  -0.05%  "
$ws.Range("E26").Value = "  -6.79%  "
$ws.Range("E27").Value = "  -0.62%  "
$ws.Range("D28").Value = "22.49"
$ws.Range("E28").Value = "  -4.85%  "
$ws.Range("E29").Value = "  +0.43%  "
$ws.Range("D30").Value = "9.05"
$ws.Range("E30").Value = "  -1.73%  "
$ws.Range("D31").Value = "148.74"
$ws.Range("E31").Value = "  -3.02%  "
$ws.Range("D32").Value = "31.57"
$ws.Range("E32").Value = "  -7.44%  "
$ws.Range("E33").Value = "  -0.19%  "
$ws.Range("D34").Value = "4.81"
$ws.Range("E34").Value = "  -7.35%  "
$ws.Range("B35").Value = "WEMIXToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D35").Value = "2.33"
$ws.Range("E35").Value = "  -3.59%  "
$ws.Range("B36").Value = "Hedera"
$ws.Range("C36").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D36").Value = "0.0690"
$ws.Range("E36").Value = "  -5.08%  "
$ws.Range("E37").Value = "  -3.87%  "
$ws.Range("D38").Value = "0.0965"
$ws.Range("E38").Value = "  -3.02%  "
$ws.Range("D39").Value = "15.09"
$ws.Range("E39").Value = "  -4.77%  "
$ws.Range("D40").Value = "2.61"
$ws.Range("D41").Value = "1.63"
$ws.Range("E41").Value = "  -3.98%  "
$ws.Range("D42").Value = "3.60"
$ws.Range("E42").Value = "  -5.31%  "
$ws.Range("D43").Value = "1.894.04"
$ws.Range("E43").Value = "  -3.69%  "
$ws.Range("D44").Value = "2.08"
$ws.Range("E44").Value = "  -12.85%  "
$ws.Range("D45").Value = "0.0257"
$ws.Range("E45").Value = "  -3.62%  "
$ws.Range("B46").Value = "FraxShare"
$ws.Range("C46").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D46").Value = "8.97"
$ws.Range("E46").Value = "  -3.02%  "
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").Value = "16.05"
$ws.Range("E47").Value = "  -9.17%  "
$ws.Range("D48").Value = "2.60"
$ws.Range("E48").Value = "  -3.30%  "
$ws.Range("D49").Value = "2.422.39"
$ws.Range("E49").Value = "  -6.89%  "
$ws.Range("D50").Value = "70.65"
$ws.Range("E50").Value = "  -0.73%  "
$ws.Range("D51").Value = "86.68"
$ws.Range("E51").Value = "  -6.78%  "

Write-Host "Applied 97 cell updates"
